$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2 - Bitcoin
Set-TextCell 2 4 "69.368.33"
$ws.Range("E2").Value = "  +0.10%  "

# Row 3 - Ethereum
Set-TextCell 3 4 "3.693.93"
$ws.Range("E3").Value = "  +0.37%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  -0.45%  "

# Row 6 - Solana
Set-TextCell 6 4 "159.39"
$ws.Range("E6").Value = "  -0.25%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.30%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.55%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -1.28%  "

# Row 11 - Cardano
Set-TextCell 11 4 "0.440"
$ws.Range("E11").Value = "  +1.36%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  -0.51%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell 13 4 "4.315.25"
$ws.Range("E13").Value = "  +0.38%  "

# Row 14 - Avalanche
Set-TextCell 14 4 "32.26"
$ws.Range("E14").Value = "  -0.97%  "

# Row 15 - WrappedEther
Set-TextCell 15 4 "3.693.41"
$ws.Range("E15").Value = "  +0.05%  "

# Row 16 - WrappedBTC
Set-TextCell 16 4 "69.397.46"
$ws.Range("E16").Value = "  +0.20%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  +3.33%  "

# Row 18 - Chainlink
Set-TextCell 18 4 "15.97"
$ws.Range("E18").Value = "  +0.72%  "

# Row 19 - Polkadot
Set-TextCell 19 4 "6.46"
$ws.Range("E19").Value = "  +0.31%  "

# Row 20 - BitcoinCash
Set-TextCell 20 4 "467.49"
$ws.Range("E20").Value = "  -0.17%  "

# Row 21 - Uniswap
Set-TextCell 21 4 "9.89"
$ws.Range("E21").Value = "  +0.32%  "

# Row 22 - Polygon
Set-TextCell 22 4 "0.651"
$ws.Range("E22").Value = "  +0.31%  "

# Row 23 - Litecoin
Set-TextCell 23 4 "80.11"
$ws.Range("E23").Value = "  +0.87%  "

# Row 24 - WrappedeETH
Set-TextCell 24 4 "3.838.98"
$ws.Range("E24").Value = "  +0.36%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.03%  "

# Row 26 - PEPE
Set-TextCell 26 4 "0.0000124"
$ws.Range("E26").Value = "  -2.69%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextCell 27 4 "10.92"
$ws.Range("E27").Value = "  -1.34%  "

# Row 28 - RenderToken
Set-TextCell 28 4 "9.14"
$ws.Range("E28").Value = "  +0.87%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +0.90%  "

# Row 30 - Fetch.AI
$ws.Range("E30").Value = "  -0.99%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  +0.31%  "

# Row 32 - NEARProtocol
Set-TextCell 32 4 "6.54"
$ws.Range("E32").Value = "  -2.39%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  -1.76%  "

# Row 34 - EthereumClassic
Set-TextCell 34 4 "26.93"
$ws.Range("E34").Value = "  +1.06%  "

# Row 35 - RenzoRestakedETH
Set-TextCell 35 4 "3.682.88"
$ws.Range("E35").Value = "  +0.90%  "

# Row 36 - Kaspa
Set-TextCell 36 4 "0.158"
$ws.Range("E36").Value = "  -1.92%  "

# Row 37 - Aptos
$ws.Range("E37").Value = "  +1.81%  "

# Row 38 - Filecoin
Set-TextCell 38 4 "6.29"
$ws.Range("E38").Value = "  +2.92%  "

# Row 39 - USDe
$ws.Range("E39").Value = "  +0.02%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  -0.05%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -1.23%  "

# Row 42 - Hedera
Set-TextCell 42 4 "0.0902"
$ws.Range("E42").Value = "  -0.19%  "

# Row 43 - Monero
Set-TextCell 43 4 "169.19"
$ws.Range("E43").Value = "  +3.42%  "

# Row 44 - Mantle
Set-TextCell 44 4 "0.941"
$ws.Range("E44").Value = "  -0.45%  "

# Row 45 - OKB
Set-TextCell 45 4 "47.12"
$ws.Range("E45").Value = "  -2.06%  "

# Row 46 - FLOKI
Set-TextCell 46 4 "0.000280"
$ws.Range("E46").Value = "  +2.24%  "

# Row 47 - dogwifhat
$ws.Range("E47").Value = "  +0.04%  "

# Row 48 - InjectiveProtocol
Set-TextCell 48 4 "27.92"
$ws.Range("E48").Value = "  -5.32%  "

# Row 49 - was SuiNetwork, now ONDO (rows 49/50 swap contents)
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextCell 49 4 "1.29"
$ws.Range("E49").Value = "  -0.43%  "

# Row 50 - was ONDO, now SuiNetwork
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextCell 50 4 "1.10"
$ws.Range("E50").Value = "  +1.17%  "

# Row 51 - Cosmos
$ws.Range("E51").Value = "  -0.94%  "
